$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("F2").Value = 4.2
$ws.Range("G2").Value = 6.4
$ws.Range("H2").Value = 1.6
$ws.Range("I2").Value = 2.04
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 7.2
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 2.16
$ws.Range("P2").Value = 1.87
$ws.Range("R2").Value = 1.38
$ws.Range("S2").Value = 2.58
$ws.Range("V2").Value = 1.96
$ws.Range("W2").Value = 1.18
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("F3").Value = 2.62
$ws.Range("G3").Value = 2.94
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3.95
$ws.Range("J3").Value = 2.62
$ws.Range("K3").Value = 3.35
$ws.Range("P3").Value = 1.41
$ws.Range("Q3").Value = 3.1
$ws.Range("L4").Value = 1.3
$ws.Range("N4").Value = 3.3
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 2.74
$ws.Range("T4").Value = 1.64
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AM5").Value = 160
$ws.Range("F5").Value = 1.66
$ws.Range("G5").Value = 1.78
$ws.Range("I5").Value = 6.8
$ws.Range("K5").Value = 4.2
$ws.Range("P5").Value = 1.78
$ws.Range("Q5").Value = 2.04
$ws.Range("U5").Value = 1.87
$ws.Range("V5").Value = 1.17
$ws.Range("W5").Value = 2.28
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("F6").Value = 1.21
$ws.Range("I6").Value = 19
$ws.Range("K6").Value = 8.800000000000001
$ws.Range("N6").Value = 6.8
$ws.Range("P6").Value = 2.76
$ws.Range("Q6").Value = 1.33
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 1.96
$ws.Range("U6").Value = 1.84
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AA7").Value = 10
$ws.Range("AD7").Value = 14.5
$ws.Range("I7").Value = 1.21
$ws.Range("J7").Value = 6.8
$ws.Range("K7").Value = 11.5
$ws.Range("N7").Value = 6.8
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 2.98
$ws.Range("Q7").Value = 1.37
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 2.08
$ws.Range("T7").Value = 2.16
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 5.5
$ws.Range("Y7").Value = 13.5
$ws.Range("F8").Value = 2.68
$ws.Range("G8").Value = 2.96
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.6
$ws.Range("M8").Value = 1.17
$ws.Range("R8").Value = 1.13
$ws.Range("AB9").Value = 7.6
$ws.Range("AC9").Value = 10
$ws.Range("AF9").Value = 11.5
$ws.Range("AG9").Value = 13
$ws.Range("G9").Value = 1.86
$ws.Range("J9").Value = 3.15
$ws.Range("N9").Value = 2.74
$ws.Range("Q9").Value = 2.42
$ws.Range("S9").Value = 4.6
$ws.Range("AA10").Value = 190
$ws.Range("AB10").Value = 7.6
$ws.Range("AC10").Value = 8.800000000000001
$ws.Range("AF10").Value = 970
$ws.Range("AG10").Value = 11
$ws.Range("AI10").Value = 120
$ws.Range("AJ10").Value = 20
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 46
$ws.Range("AM10").Value = 180
$ws.Range("AN10").Value = 970
$ws.Range("F10").Value = 1.74
$ws.Range("G10").Value = 1.88
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 3.45
$ws.Range("K10").Value = 3.95
$ws.Range("M10").Value = 1.09
$ws.Range("N10").Value = 3.05
$ws.Range("O10").Value = 1.39
$ws.Range("Q10").Value = 2.16
$ws.Range("R10").Value = 1.26
$ws.Range("S10").Value = 4
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 1.82
$ws.Range("V10").Value = 1.19
$ws.Range("W10").Value = 2.12
$ws.Range("X10").Value = 1000
$ws.Range("Z10").Value = 46
$ws.Range("AA11").Value = 75
$ws.Range("AB11").Value = 10
$ws.Range("AD11").Value = 17
$ws.Range("AE11").Value = 60
$ws.Range("AF11").Value = 22
$ws.Range("AG11").Value = 14.5
$ws.Range("AH11").Value = 28
$ws.Range("AK11").Value = 48
$ws.Range("AL11").Value = 75
$ws.Range("AN11").Value = 60
$ws.Range("AO11").Value = 65
$ws.Range("P11").Value = 1.49
$ws.Range("Q11").Value = 2.66
$ws.Range("X11").Value = 8.6
$ws.Range("Y11").Value = 9.6
$ws.Range("Z11").Value = 22
$ws.Range("AA12").Value = 75
$ws.Range("AB12").Value = 9.4
$ws.Range("AC12").Value = 8
$ws.Range("AD12").Value = 970
$ws.Range("AE12").Value = 48
$ws.Range("AF12").Value = 15
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 970
$ws.Range("AI12").Value = 60
$ws.Range("AJ12").Value = 32
$ws.Range("AK12").Value = 27
$ws.Range("AL12").Value = 44
$ws.Range("AM12").Value = 130
$ws.Range("AN12").Value = 22
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 3.4
$ws.Range("O12").Value = 1.36
$ws.Range("Q12").Value = 2.06
$ws.Range("R12").Value = 1.3
$ws.Range("S12").Value = 3.75
$ws.Range("T12").Value = 1.79
$ws.Range("U12").Value = 2.06
$ws.Range("W12").Value = 1.75
$ws.Range("Y12").Value = 13.5
$ws.Range("Z12").Value = 27
$ws.Range("AB13").Value = 9.4
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 21
$ws.Range("AE13").Value = 70
$ws.Range("AF13").Value = 16
$ws.Range("AG13").Value = 14
$ws.Range("AH13").Value = 26
$ws.Range("AI13").Value = 90
$ws.Range("AJ13").Value = 40
$ws.Range("AK13").Value = 36
$ws.Range("AL13").Value = 65
$ws.Range("AN13").Value = 34
$ws.Range("F13").Value = 2.16
$ws.Range("G13").Value = 2.64
$ws.Range("H13").Value = 3.65
$ws.Range("I13").Value = 4.5
$ws.Range("J13").Value = 2.74
$ws.Range("K13").Value = 3.5
$ws.Range("M13").Value = 1.1
$ws.Range("N13").Value = 2.82
$ws.Range("O13").Value = 1.44
$ws.Range("R13").Value = 1.23
$ws.Range("S13").Value = 4.1
$ws.Range("T13").Value = 1.96
$ws.Range("U13").Value = 1.84
$ws.Range("V13").Value = 1.29
$ws.Range("W13").Value = 1.7
$ws.Range("X13").Value = 12.5
$ws.Range("Y13").Value = 14.5
$ws.Range("Z13").Value = 34
$ws.Range("AA14").Value = 190
$ws.Range("AC14").Value = 8.6
$ws.Range("AE14").Value = 100
$ws.Range("AF14").Value = 9.4
$ws.Range("AN14").Value = 11.5
$ws.Range("AO14").Value = 130
$ws.Range("G14").Value = 1.71
$ws.Range("J14").Value = 3.9
$ws.Range("N14").Value = 3.55
$ws.Range("T14").Value = 1.99
$ws.Range("U14").Value = 1.89
$ws.Range("W14").Value = 2.4
